# Update to new UI
# The workbook stores three JDS (job-data-set) sheets, each with a timeZoneId
# value in cell B5. The old raw Windows time-zone id strings are being
# replaced with the new UI's display strings (e.g. "(UTC-10:00) Hawaii").

$wb = $excel.ActiveWorkbook

# --- Sheet "initial_JDS" ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("initial_JDS")
$ws1.Activate()
$ws1.Range("B5").Value = "(UTC-10:00) Hawaii"
$ws1.Range("B5").Select()

# --- Sheet "updated_JDS" ----------------------------------------------------
$ws2 = $wb.Worksheets.Item("updated_JDS")
$ws2.Activate()
$ws2.Range("B5").Value = "(UTC+10:00) Hobart"
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
$ws2.Range("B5").Select()

# --- Sheet "nozip_JDS" -------------------------------------------------------
$ws3 = $wb.Worksheets.Item("nozip_JDS")
$ws3.Activate()
$ws3.Range("B5").Value = "(UTC-08:00) Pacific Time (US & Canada)"
$ws3.Columns("B").ColumnWidth = 33
$ws3.Columns("B").Select()
